# The published dataset gained one more weekly record for "Ají" (Hortaliza,
# Vega Central Mapocho de Santiago). It was inserted as a new row 146,
# pushing every following row down by one (old row 146 -> new row 147, ...,
# old row 220 -> new row 221), and the sheet's used range grew from
# A1:R220 to A1:R221.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 146; Excel shifts rows 146..220 down to
# 147..221 automatically and grows the sheet dimension for us.
$ws.Rows.Item(146).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A146").Value = 9
$ws.Range("B146").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C146").Value = "Metropolitana"
$ws.Range("D146").Value = 44572
$ws.Range("E146").Value = 13
$ws.Range("F146").Value = 100112021
$ws.Range("G146").Value = "Ají"
$ws.Range("H146").Value = "Americana (o)"
$ws.Range("I146").Value = "Primera"
$ws.Range("J146").Value = 16
$ws.Range("K146").Value = 26000
$ws.Range("L146").Value = 28000
$ws.Range("M146").Value = 27000
$ws.Range("N146").Value = '$/caja 25 kilos'
$ws.Range("O146").Value = "Región Metropolitana"
$ws.Range("P146").Value = 1080
$ws.Range("Q146").Value = 25
$ws.Range("R146").Value = "Hortaliza"
